$d = $word.ActiveDocument

# First paragraph holds the document's ID placeholder marker.
$p1 = $d.Paragraphs(1)

# Add a (line-less) paragraph border with 5pt spacing on every edge,
# matching the "highlighted id" paragraph style already used elsewhere
# in this template (see the (b)(2) paragraph further down).
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromRight = 5

# Bump the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.LeftIndent = 11.25

# Update the ID placeholder text and drop the trailing space run that
# used to follow it.
$d.Content.Find.Execute("**ID__AFFARS_pgi_5316_topic_12__ID** ", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SMC_PGI_5316_505__ID**", 2)
